# Update stack trace line numbers / identifiers to reflect the move from
# POI 3.17.0 to 4.0.1, and replace the Maven/Tycho/Equinox launch frames
# with the Eclipse JDT JUnit runner frames.

$d = $word.ActiveDocument

function Replace-Literal($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Simple line-number / identifier swaps.
Replace-Literal "JavaMethodService.java:163" "JavaMethodService.java:162"
Replace-Literal "AbstractService.java:136" "AbstractService.java:135"
Replace-Literal "EvaluationServices.java:168" "EvaluationServices.java:172"
Replace-Literal "AstSwitch.java:118" "AstSwitch.java:119"
Replace-Literal "AbstractTemplatesTestSuite.java:480" "AbstractTemplatesTestSuite.java:462"
Replace-Literal "AbstractTemplatesTestSuite.java:389" "AbstractTemplatesTestSuite.java:372"
Replace-Literal "GeneratedMethodAccessor74" "GeneratedMethodAccessor75"

# Replace the Maven Surefire / Tycho / Equinox launcher frames with the
# Eclipse JDT JUnit runner frames. This block is far longer than the
# 255-character limit Word imposes on Find/Replace text, so locate its
# first and last (unique) lines and overwrite the Range between them
# directly instead of going through Find.Execute's ReplaceWith.
$rStart = $d.Content.Duplicate
$rStart.Find.Execute("`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)", `
                      $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$rEnd = $d.Content.Duplicate
$rEnd.Find.Execute("at org.eclipse.equinox.launcher.Main.main(Main.java:1471)", `
                    $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$block = $d.Range($rStart.Start, $rEnd.End)

$newBlock = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n" + `
            "`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n" + `
            "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n" + `
            "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n" + `
            "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n" + `
            "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"

$block.Text = $newBlock
